$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "devyani"
$ws.Range("B2").Value = "devyanikumar947@gmail.com"
$ws.Range("C2").Value = "Devyani@123"
